# StaffInfo 오타 수정
# Fix typo: the "type" row (row 2) for the "Staff" column incorrectly
# used the shared string "String" (capitalized) while every other type
# cell uses lowercase type names (int, bool). Correct it to "string".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "string"

# Leave the selection where the editor ended up after making the fix.
$ws.Range("I8").Select()
